$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# Cells whose new text would otherwise be auto-parsed as a plain number (e.g. "565.33")
# get a leading apostrophe so Excel keeps them as literal text, matching the source data
# (mixed "59.315.55"-style thousand-grouped values live in the same column), then we strip
# the resulting cell formatting so no stray style is left behind.

$ws.Range("D2").Value = "59.315.55"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.982.09"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'565.33"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "'137.05"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.517"
$ws.Range("D9").Value = "2.981.63"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'0.133"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "'5.27"
$ws.Range("E11").Value = "  +7.83%  "
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "'33.56"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "3.474.02"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "2.978.98"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "59.304.09"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "'435.21"
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("D21").Value = "'13.67"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").Value = "'7.03"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").Value = "'13.00"
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("D25").Value = "'80.02"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  +5.71%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'7.72"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'25.68"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("D36").Value = "0.0₃0763"
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("D37").Value = "'2.05"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").Value = "'48.62"
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "'2.81"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").Value = "'395.09"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "2.710.31"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'122.70"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D48").Value = "'34.50"
$ws.Range("E48").Value = "  +10.48%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("D51").Value = "'23.13"

# Strip the auto-applied number formatting / quote-prefix marker so the cells look
# exactly like the rest of the untouched text cells in the sheet.
$ws.Range("D5,D6,D8,D10,D11,D14,D20,D21,D23,D24,D25,D30,D31,D32,D37,D38,D40,D41,D47,D48,D51").ClearFormats()
